$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new learning-log rows (row 9 and row 10) at the bottom of the
# existing table, matching the pattern used by the prior rows.
$ws.Cells.Item(9, 1).Value = 25
$ws.Cells.Item(9, 2).Value = "1：06—3：53"
$ws.Cells.Item(9, 3).Value = "数据类型，变量，常量"

$ws.Cells.Item(10, 1).Value = 26
$ws.Cells.Item(10, 2).Value = "3：06—6：53"
$ws.Cells.Item(10, 3).Value = "常量 字符串 转义字符 for while 函数与数组"

# Move the active selection to the last edited cell, like Excel would after
# typing the final entry.
$ws.Range("C10").Select()
